# Auto-generated Excel COM-interop script to apply the diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 32328
$ws.Range("J81").Value = 32328
$ws.Range("L81").Value = 32328
$ws.Range("N81").Value = -34324
$ws.Range("H84").Value = 32328
$ws.Range("J84").Value = 32328
$ws.Range("L84").Value = 96984
$ws.Range("N84").Value = -106968
$ws.Range("H111").Value = 3044.0557
$ws.Range("I111").Value = 1138.3077
$ws.Range("J111").Value = 7999
$ws.Range("K111").Value = 3414.9231
$ws.Range("L111").Value = 23997
$ws.Range("M111").Value = -347.9231
$ws.Range("N111").Value = -30131
$ws.Range("H112").Value = 2942.6667
$ws.Range("J112").Value = 2975.1724
$ws.Range("L112").Value = 8925.5172
$ws.Range("N112").Value = -11141.5172
$ws.Range("H129").Value = 821.1
$ws.Range("I129").Value = 608.375
$ws.Range("J129").Value = 898.4545000000001
$ws.Range("K129").Value = 1825.125
$ws.Range("L129").Value = 2695.3635
$ws.Range("M129").Value = 3174.875
$ws.Range("N129").Value = -12695.3635
$ws.Range("H132").Value = 2432.1562
$ws.Range("I132").Value = 2432.1562
$ws.Range("K132").Value = 7296.4686
$ws.Range("M132").Value = -4766.4686
$ws.Range("H135").Value = 50015052
$ws.Range("I135").Value = 1084
$ws.Range("K135").Value = 9756
$ws.Range("M135").Value = -7221
$ws.Range("H137").Value = 265591.3
$ws.Range("I137").Value = 574120.1
$ws.Range("J137").Value = 85616.164
$ws.Range("K137").Value = 1722360.3
$ws.Range("L137").Value = 256848.492
$ws.Range("M137").Value = -1719810.3
$ws.Range("N137").Value = -261948.492
$ws.Range("H138").Value = 1704.9375
$ws.Range("I138").Value = 529.4828
$ws.Range("K138").Value = 1588.4484
$ws.Range("M138").Value = 3551.5516

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2136.1765
$ws.Range("I2").Value = 1275.25
$ws.Range("J2").Value = 4202.4
$ws.Range("K2").Value = 1275.25
$ws.Range("L2").Value = 4202.4
$ws.Range("M2").Value = -1162.25
$ws.Range("N2").Value = -4428.4
$ws.Range("H32").Value = 18270.709
$ws.Range("I32").Value = 19977.291
$ws.Range("J32").Value = 4861.857
$ws.Range("K32").Value = 19977.291
$ws.Range("L32").Value = 4861.857
$ws.Range("M32").Value = -19690.291
$ws.Range("N32").Value = -5435.857
$ws.Range("H116").Value = 2136.1765
$ws.Range("I116").Value = 1275.25
$ws.Range("J116").Value = 4202.4
$ws.Range("K116").Value = 1275.25
$ws.Range("L116").Value = 4202.4
$ws.Range("M116").Value = 1018.75
$ws.Range("N116").Value = -8790.4
$ws.Range("H132").Value = 16015.306
$ws.Range("I132").Value = 1642.1578
$ws.Range("J132").Value = 32079.412
$ws.Range("K132").Value = 4926.4734
$ws.Range("L132").Value = 96238.236
$ws.Range("M132").Value = -2396.4734
$ws.Range("N132").Value = -101298.236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2136.1765
$ws.Range("I3").Value = 1275.25
$ws.Range("J3").Value = 4202.4
$ws.Range("K3").Value = 1275.25
$ws.Range("L3").Value = 4202.4
$ws.Range("M3").Value = -1161.25
$ws.Range("N3").Value = -4430.4
$ws.Range("H105").Value = 2779936.5
$ws.Range("I105").Value = 2208
$ws.Range("J105").Value = 4547582
$ws.Range("K105").Value = 2208
$ws.Range("L105").Value = 4547582
$ws.Range("M105").Value = -461
$ws.Range("N105").Value = -4551076
$ws.Range("H107").Value = 1831.8572
$ws.Range("I107").Value = 950
$ws.Range("J107").Value = 3007.6667
$ws.Range("K107").Value = 950
$ws.Range("L107").Value = 3007.6667
$ws.Range("M107").Value = 970
$ws.Range("N107").Value = -6847.6667
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 49944.25
$ws.Range("J129").Value = 49944.25
$ws.Range("L129").Value = 49944.25
$ws.Range("N129").Value = -59944.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1158.2632
$ws.Range("I16").Value = 927.9091
$ws.Range("K16").Value = 927.9091
$ws.Range("M16").Value = -640.9091
$ws.Range("H31").Value = 12988.921
$ws.Range("I31").Value = 21047.316
$ws.Range("J31").Value = 4930.5264
$ws.Range("K31").Value = 21047.316
$ws.Range("L31").Value = 4930.5264
$ws.Range("M31").Value = -20752.316
$ws.Range("N31").Value = -5520.5264
$ws.Range("H34").Value = 12988.921
$ws.Range("I34").Value = 21047.316
$ws.Range("J34").Value = 4930.5264
$ws.Range("K34").Value = 21047.316
$ws.Range("L34").Value = 4930.5264
$ws.Range("M34").Value = -20845.316
$ws.Range("N34").Value = -5334.5264
$ws.Range("H58").Value = 32779.938
$ws.Range("I58").Value = 1842.3
$ws.Range("J58").Value = 84342.664
$ws.Range("K58").Value = 1842.3
$ws.Range("L58").Value = 84342.664
$ws.Range("M58").Value = -1639.3
$ws.Range("N58").Value = -84748.664
$ws.Range("H99").Value = 4422
$ws.Range("I99").Value = 3290.8
$ws.Range("J99").Value = 7250
$ws.Range("K99").Value = 3290.8
$ws.Range("L99").Value = 7250
$ws.Range("M99").Value = -1792.8
$ws.Range("N99").Value = -10246
$ws.Range("H100").Value = 300000
$ws.Range("J100").Value = 300000
$ws.Range("L100").Value = 300000
$ws.Range("N100").Value = -302164
$ws.Range("H113").Value = 1158.2632
$ws.Range("I113").Value = 927.9091
$ws.Range("K113").Value = 927.9091
$ws.Range("M113").Value = 1242.0909
$ws.Range("H126").Value = 4422
$ws.Range("I126").Value = 3290.8
$ws.Range("J126").Value = 7250
$ws.Range("K126").Value = 9872.400000000001
$ws.Range("L126").Value = 21750
$ws.Range("M126").Value = -7402.400000000001
$ws.Range("N126").Value = -26690
$ws.Range("H134").Value = 860.65515
$ws.Range("I134").Value = 658.36
$ws.Range("J134").Value = 2125
$ws.Range("K134").Value = 1975.08
$ws.Range("L134").Value = 6375
$ws.Range("M134").Value = 559.9200000000001
$ws.Range("N134").Value = -11445
$ws.Range("H136").Value = 32779.938
$ws.Range("I136").Value = 1842.3
$ws.Range("J136").Value = 84342.664
$ws.Range("K136").Value = 5526.9
$ws.Range("L136").Value = 253027.992
$ws.Range("M136").Value = -2976.9
$ws.Range("N136").Value = -258127.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1130.5625
$ws.Range("I5").Value = 934.5714
$ws.Range("K5").Value = 2803.7142
$ws.Range("M5").Value = -2691.7142
$ws.Range("H62").Value = 8238.429
$ws.Range("J62").Value = 8238.429
$ws.Range("L62").Value = 24715.287
$ws.Range("N62").Value = -26087.287
$ws.Range("H65").Value = 8238.429
$ws.Range("J65").Value = 8238.429
$ws.Range("L65").Value = 74145.861
$ws.Range("N65").Value = -81009.861
$ws.Range("H131").Value = 752.58
$ws.Range("J131").Value = 752.58
$ws.Range("L131").Value = 2257.74
$ws.Range("N131").Value = -12337.74
$ws.Range("H135").Value = 1130.5625
$ws.Range("I135").Value = 934.5714
$ws.Range("K135").Value = 8411.142600000001
$ws.Range("M135").Value = -5876.142600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 56530.855
$ws.Range("I132").Value = 52892.4
$ws.Range("K132").Value = 158677.2
$ws.Range("M132").Value = -156147.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 380.06668
$ws.Range("I16").Value = 375.1
$ws.Range("J16").Value = 390
$ws.Range("K16").Value = 375.1
$ws.Range("L16").Value = 390
$ws.Range("M16").Value = -205.1
$ws.Range("N16").Value = -730
$ws.Range("H46").Value = 1052.8572
$ws.Range("I46").Value = 826.25
$ws.Range("J46").Value = 1778
$ws.Range("K46").Value = 826.25
$ws.Range("L46").Value = 1778
$ws.Range("M46").Value = -638.25
$ws.Range("N46").Value = -2154
$ws.Range("H132").Value = 465649.7
$ws.Range("I132").Value = 636131
$ws.Range("J132").Value = 2914.7144
$ws.Range("K132").Value = 1908393
$ws.Range("L132").Value = 8744.143199999999
$ws.Range("M132").Value = -1905863
$ws.Range("N132").Value = -13804.1432
$ws.Range("H136").Value = 45166.668
$ws.Range("I136").Value = 64875
$ws.Range("K136").Value = 194625
$ws.Range("M136").Value = -192075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 760.4761999999999
$ws.Range("I126").Value = 660.6
$ws.Range("K126").Value = 1981.8
$ws.Range("M126").Value = 488.1999999999998
$ws.Range("H132").Value = 2367.1667
$ws.Range("I132").Value = 2031.8334
$ws.Range("K132").Value = 6095.5002
$ws.Range("M132").Value = -3565.5002
